$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value = "6.5_19.5_340km"
$ws.Range("B2").Value = 420.0
$ws.Range("C2").Value = 1140.0
$ws.Range("D2").Value = 330.0
$ws.Range("F2").Value = 2640.0
$ws.Range("G2").Value = 323.0
$ws.Range("I2").Value = 4080.0
$ws.Range("J2").Value = 355.0
$ws.Range("L2").Value = 5460.0
$ws.Range("M2").Value = 335.0
$ws.Range("O2").Value = 6960.0
$ws.Range("P2").Value = 327.0
$ws.Range("R2").Value = 8400.0
$ws.Range("S2").Value = 328.0
$ws.Range("T2").Value = 9000.0
$ws.Range("U2").Value = 9840.0
$ws.Range("V2").Value = 350.0

$ws.Range("A11").Value = "6.5_18.25_280km"
$ws.Range("B11").Value = 360.0
$ws.Range("D11").Value = 291.0
$ws.Range("E11").Value = 1800.0
$ws.Range("F11").Value = 2520.0
$ws.Range("G11").Value = 282.0
$ws.Range("H11").Value = 3300.0
$ws.Range("J11").Value = 285.0
$ws.Range("L11").Value = 5400.0
$ws.Range("M11").Value = 265.0
$ws.Range("N11").Value = 6120.0
$ws.Range("P11").Value = 301.0
$ws.Range("R11").Value = 8280.0
$ws.Range("S11").Value = 268.0
$ws.Range("T11").Value = 9000.0
$ws.Range("V11").Value = 260.0

$ws.Range("A12").Value = "6_18_250km"
$ws.Range("D12").Value = 230.0
$ws.Range("E12").Value = 1800.0
$ws.Range("G12").Value = 235.0
$ws.Range("H12").Value = 3240.0
$ws.Range("I12").Value = 3960.0
$ws.Range("J12").Value = 241.0
$ws.Range("M12").Value = 276.0
$ws.Range("P12").Value = 250.0
$ws.Range("S12").Value = 240.0
